$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New row 7: finalization of MLD/MCD, entered first so its string gets the
# next shared-string slot (matches the order strings were introduced).
$ws.Range("B7").Value = "Finalisation du MLD/MCD."

# Row 5: hours corrected from "5 heures " to "4.5 heures "
$ws.Range("C5").Value = "4.5 heures "

# New row 6: half hour entry
$ws.Range("C6").Value = "0.5 heure"

# New row 6: start of documentation work
$ws.Range("B6").Value = "Commencement de la documentation du projet. Rédiger l'introduction, les objectifs et la planification initiale."

# Dates for the two new rows - copy the date formatting from A5 so the same
# cell style (numFmtId 14) is reused instead of minting a new style.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A6").Value = "05/09/2018"
$ws.Range("A7").Value = "05/15/2018"

# Row 7 hours
$ws.Range("C7").Value = "1 heure"

$ws.PageSetup.Orientation = 2

$ws.Range("B6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120
